$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A7 is an empty text cell in the source data (mirrors A2/A3/A5/A6).
# A bare "'" forces Excel to store it as an empty *text* value rather than
# clearing the cell outright; resetting the style afterwards drops the
# quote-prefix formatting that the apostrophe trick leaves behind.
$ws.Range("A7").Value = "'"
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = "احمد"

# C7 ("22") must round-trip as text (numberStoredAsText), matching the
# other numeric-looking text values already in the sheet (e.g. C6).
# NumberFormat="@" forces text storage; resetting the style afterwards
# drops that explicit formatting so the cell keeps the sheet's default style.
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "22"
$ws.Range("C7").Style = "Normal"

$ws.Range("D7").Value = "الجزائري"
$ws.Range("E7").Value = "الرحلة 1"
$ws.Range("F7").Value = "C3"
$ws.Range("G7").Value = "NRC"
$ws.Range("H7").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٣٨:١١ م"
